$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, whether it is the numeric "Price"
# column (which must be forced to text so Excel does not silently
# reinterpret strings like "1.001" or "0.9998" as numbers).
$updates = @(
    @{ Cell = 'D2'; Value = '30.253.86'; ForceText = $true }
    @{ Cell = 'E2'; Value = '  -0.25%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '1.913.66'; ForceText = $true }
    @{ Cell = 'E3'; Value = '  -1.02%  '; ForceText = $false }
    @{ Cell = 'D4'; Value = '0.9998'; ForceText = $true }
    @{ Cell = 'E4'; Value = '  -0.07%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '0.7415'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  -0.87%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '244.07'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  -1.89%  '; ForceText = $false }
    @{ Cell = 'E7'; Value = '  +0.03%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '0.3148'; ForceText = $true }
    @{ Cell = 'E8'; Value = '  -1.90%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '27.13'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  -4.18%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '0.06973'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  -1.93%  '; ForceText = $false }
    @{ Cell = 'D11'; Value = '0.7852'; ForceText = $true }
    @{ Cell = 'E11'; Value = '  -0.50%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '0.07977'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  -0.30%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '1.934.53'; ForceText = $true }
    @{ Cell = 'E13'; Value = '  +0.04%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '5.285'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  -1.90%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '91.71'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  -3.03%  '; ForceText = $false }
    @{ Cell = 'B16'; Value = 'Avalanche'; ForceText = $false }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; ForceText = $false }
    @{ Cell = 'D16'; Value = '14.31'; ForceText = $true }
    @{ Cell = 'E16'; Value = '  -2.29%  '; ForceText = $false }
    @{ Cell = 'B17'; Value = 'WrappedBTC'; ForceText = $false }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; ForceText = $false }
    @{ Cell = 'D17'; Value = '30.317.93'; ForceText = $true }
    @{ Cell = 'E17'; Value = '  -0.04%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '245.31'; ForceText = $true }
    @{ Cell = 'E18'; Value = '  -3.14%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '5.829'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  +0.49%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '0.000007836'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  -2.52%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '2.223.00'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  +1.77%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '1.001'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  +0.01%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '0.9993'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  -0.07%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '6.665'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  -2.23%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '9.439'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  -1.49%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '165.28'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  +0.46%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '19.05'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  -0.37%  '; ForceText = $false }
    @{ Cell = 'D28'; Value = '0.1273'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  -4.35%  '; ForceText = $false }
    @{ Cell = 'D29'; Value = '2.136'; ForceText = $true }
    @{ Cell = 'E29'; Value = '  -8.82%  '; ForceText = $false }
    @{ Cell = 'D30'; Value = '1.351'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  -0.43%  '; ForceText = $false }
    @{ Cell = 'D31'; Value = '1.549'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  +1.22%  '; ForceText = $false }
    @{ Cell = 'D32'; Value = '4.329'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  -2.70%  '; ForceText = $false }
    @{ Cell = 'D33'; Value = '4.092'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  -1.39%  '; ForceText = $false }
    @{ Cell = 'D34'; Value = '0.05224'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  +1.66%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '1.299'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  +1.14%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '0.7542'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  +0.58%  '; ForceText = $false }
    @{ Cell = 'D37'; Value = '2.757'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  -0.33%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '0.01945'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  -1.20%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '2.786'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  -0.66%  '; ForceText = $false }
    @{ Cell = 'D40'; Value = '6.426'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  +0.11%  '; ForceText = $false }
    @{ Cell = 'D41'; Value = '76.27'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  -2.50%  '; ForceText = $false }
    @{ Cell = 'D42'; Value = '0.4507'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  -0.27%  '; ForceText = $false }
    @{ Cell = 'D43'; Value = '1.955'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  -1.99%  '; ForceText = $false }
    @{ Cell = 'E44'; Value = '  -0.03%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '0.8320'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  -1.36%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '7.725'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  +2.47%  '; ForceText = $false }
    @{ Cell = 'B47'; Value = 'Quant'; ForceText = $false }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; ForceText = $false }
    @{ Cell = 'D47'; Value = '101.55'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  -1.01%  '; ForceText = $false }
    @{ Cell = 'B48'; Value = 'EnergySwap'; ForceText = $false }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; ForceText = $false }
    @{ Cell = 'D48'; Value = '9.900'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  +0.70%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '2.107.65'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  +0.92%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '37.14'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  -0.96%  '; ForceText = $false }
    @{ Cell = 'D51'; Value = '0.1219'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  +1.34%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Preserve the original "General" style: set a text format only
        # long enough to force string storage, then clear the format again
        # so no stray numFmt/style survives on the cell.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.ClearFormats()
    } else {
        $rng.Value = $u.Value
    }
}
